$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.249.24"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "3.604.33"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'604.67"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").Value = "'195.55"
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("D9").Value = "'0.206"
$ws.Range("E9").Value = "  -2.03%  "
$ws.Range("E10").Value = "  -1.80%  "
$ws.Range("D11").Value = "'53.71"
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("D12").Value = "'0.0000303"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").Value = "4.175.48"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").Value = "'13.03"
$ws.Range("E15").Value = "  +2.48%  "
$ws.Range("D16").Value = "'595.97"
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("D17").Value = "70.393.47"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.600.07"
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'19.02"
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").Value = "'17.79"
$ws.Range("E22").Value = "  -2.91%  "
$ws.Range("D23").Value = "'5.17"
$ws.Range("E23").Value = "  -2.48%  "
$ws.Range("D24").Value = "'102.19"
$ws.Range("E24").Value = "  -1.17%  "
$ws.Range("D25").Value = "'4.61"
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("E26").Value = "  -3.04%  "
$ws.Range("D27").Value = "'10.74"
$ws.Range("E27").Value = "  -2.16%  "
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("D30").Value = "'4.76"
$ws.Range("E30").Value = "  +4.86%  "
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("E32").Value = "  -3.79%  "
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("D34").Value = "'63.19"
$ws.Range("E34").Value = "  -0.46%  "
$ws.Range("D35").Value = "0.0₃0898"
$ws.Range("E35").Value = "  +6.85%  "
$ws.Range("D36").Value = "3.897.52"
$ws.Range("E36").Value = "  +4.22%  "
$ws.Range("D37").Value = "'530.51"
$ws.Range("E37").Value = "  +7.07%  "
$ws.Range("D38").Value = "'3.11"
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("D41").Value = "'0.390"
$ws.Range("E41").Value = "  -1.50%  "
$ws.Range("E42").Value = "  -3.10%  "
$ws.Range("E43").Value = "  -2.21%  "
$ws.Range("D44").Value = "'0.0453"
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("D45").Value = "'3.41"
$ws.Range("E45").Value = "  +2.39%  "
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").Value = "'8.61"
$ws.Range("E48").Value = "  -0.98%  "
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("D50").Value = "'0.000250"
$ws.Range("E50").Value = "  +1.42%  "
$ws.Range("E51").Value = "  +0.39%  "
